$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the Kyrgyz NEET indicator title text in A1 (wording tweak + trimmed leading spaces)
$ws.Range("A1").Value = "8.6.1 Иштебеген, окубаган жана кесиптик көндүмдөрдү үйрөнбөгөн (15 жаштан 24 жашка чейинки ) жаштардын үлүшү  "

# 2. Row 1 custom height changes from 54 to 48
$ws.Rows(1).RowHeight = 48

# 3. A1 vertical alignment changes from Top to Center (text stays left/wrap)
$ws.Range("A1").VerticalAlignment = -4108   # xlCenter

# 4. Add the new 2023 column (T) of data, matching the formatting of column S (2022)
$ws.Range("S4:S7").Copy($ws.Range("T4:T7"))
$ws.Range("T4").Value = 2023
$ws.Range("T5").Value = 18.6
$ws.Range("T6").Value = 11.5
$ws.Range("T7").Value = 25.9
$excel.CutCopyMode = 0

# 5. Reset the active selection back to A1 (workbook was saved with focus on A1, not S4:S7)
$null = $ws.Range("A1").Select()
